# DAF Core creation - device, goal, careteam
# Populate the "Dev Complete" column (D) on the Names-key sheet to reflect
# progress on newly-created/updated profiles, and update the Owner note
# for the Conformance row from "Both..." to "Discuss".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Names-key")

# Conformance row: Owner column changes from "Both..." to "Discuss"
$ws.Range("C5").Value = "Discuss"

# Dev Complete column (D) updates
$ws.Range("D2").Value = "Initial"
$ws.Range("D3").Value = "Initial"
$ws.Range("D4").Value = "Created, not hooked in"
$ws.Range("D8").Value = "created"
$ws.Range("D9").Value = "created"
$ws.Range("D10").Value = "created"
$ws.Range("D11").Value = "created"
$ws.Range("D12").Value = "created"
$ws.Range("D13").Value = "created"
$ws.Range("D14").Value = "created"
$ws.Range("D15").Value = "created"
$ws.Range("D16").Value = "created"
$ws.Range("D17").Value = "created"
$ws.Range("D18").Value = "created"
$ws.Range("D19").Value = "created"
$ws.Range("D20").Value = "Initial"
$ws.Range("D21").Value = "created"
$ws.Range("D22").Value = "created"
$ws.Range("D23").Value = "created"
$ws.Range("D24").Value = "created"
$ws.Range("D25").Value = "created"

# Widen column D to fit the new text, and move the active selection
# (mirrors the author's last on-screen position when saving).
$ws.Columns.Item(4).ColumnWidth = 20.6
$ws.Range("D7").Select()
